$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mutual information values per document (column B) updated for rows below.
$values = @{
    2 = 5
    4 = 1
    5 = 2
    8 = 4
    9 = 2
    11 = 5
    13 = 4
    15 = 2
    16 = 5
    17 = 3
    18 = 2
    20 = 2
    21 = 4
    22 = 1
    23 = 4
    24 = 5
    25 = 4
    26 = 1
    27 = 2
    28 = 1
    29 = 5
    30 = 1
    31 = 3
    32 = 2
    33 = 2
    34 = 4
    35 = 2
    36 = 3
    37 = 5
    38 = 1
    39 = 1
    40 = 4
    42 = 5
    43 = 2
    44 = 1
    45 = 1
    48 = 4
    50 = 2
    51 = 1
    52 = 4
    54 = 1
    55 = 5
    56 = 1
    57 = 2
    58 = 5
    59 = 5
    60 = 5
    61 = 3
    62 = 5
    63 = 2
    65 = 2
    66 = 4
    69 = 2
    72 = 5
    74 = 3
    75 = 1
    77 = 3
    78 = 4
    79 = 1
    80 = 4
    81 = 2
    82 = 5
    83 = 5
    84 = 5
    85 = 2
    86 = 5
    87 = 2
    88 = 1
    89 = 1
    90 = 5
    91 = 2
}

foreach ($row in $values.Keys) {
    $ws.Cells.Item([int]$row, 2).Value = $values[$row]
}

